# Auto-generated Excel COM-interop script to update cryptos list (ranking refresh).
# For each row whose market data changed: Price (D) and Volume(1h) (E) are updated.
# Rows whose ranking position shifted also get a new Coin (B) and Link (C).
# D/E columns are forced to Text format so numeric-looking strings (e.g. "0.9994")
# are stored as text, matching the original inlineStr cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.468.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.54%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.729.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.62%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4799"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.71%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2671"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06231"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.30%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.47"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07127"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.73"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.26%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.554"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.22%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.17"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.10%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.473.95"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.58%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006953"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.61%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.958.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.562"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.928"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.319"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.50"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.74%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.34"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.794"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.409"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.32%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.75"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.992"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08007"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.740"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04562"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.40%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.616"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.73%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6432"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.77%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9919"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.97%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9421"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.21%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.998"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.83%  "

# Row 39
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "107.73"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.31%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.407"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.76%  "

# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.74%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01502"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.99%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.652"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.87%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3910"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.82%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.958"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +12.89%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1194"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.51%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05324"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.37%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.82"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.25%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.866"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.01%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.273"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.91%  "

# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3430"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.25%  "

